# Update the "Recorded By" column (G) values for specific rows: reorder the
# comma-separated list of recorder names/emails (same members, new order) to
# match the synced upstream data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = "system, System, backup@backdoor.com"
$ws.Range("G3").Value = "System, dnasr281@gmail.com"
$ws.Range("G4").Value = "System, backup@backdoor.com"
$ws.Range("G5").Value = "System, backup@backdoor.com"
$ws.Range("G6").Value = "System, dnasr281@gmail.com"
$ws.Range("G7").Value = "admin@admin.com, System"
$ws.Range("G8").Value = "System, backup@backdoor.com"

$ws.Range("G28").Value = "system, System, backup@backdoor.com"
$ws.Range("G29").Value = "System, dnasr281@gmail.com"
$ws.Range("G30").Value = "System, backup@backdoor.com"
$ws.Range("G31").Value = "System, backup@backdoor.com"
$ws.Range("G32").Value = "System, dnasr281@gmail.com"
$ws.Range("G33").Value = "admin@admin.com, System"
$ws.Range("G34").Value = "System, backup@backdoor.com"

$ws.Range("G54").Value = "system, System, backup@backdoor.com"
$ws.Range("G55").Value = "System, dnasr281@gmail.com"
$ws.Range("G56").Value = "System, backup@backdoor.com"
$ws.Range("G57").Value = "System, backup@backdoor.com"
$ws.Range("G58").Value = "System, dnasr281@gmail.com"
$ws.Range("G59").Value = "admin@admin.com, System"
$ws.Range("G60").Value = "System, backup@backdoor.com"

$ws.Range("G80").Value = "System, backup@backdoor.com"
$ws.Range("G81").Value = "System, backup@backdoor.com"
$ws.Range("G82").Value = "System, backup@backdoor.com"

$ws.Range("G87").Value = "admin@admin.com, dnasr281@gmail.com"

$ws.Range("G106").Value = "System, backup@backdoor.com"
$ws.Range("G107").Value = "System, backup@backdoor.com"
$ws.Range("G108").Value = "System, backup@backdoor.com"

$ws.Range("G113").Value = "admin@admin.com, dnasr281@gmail.com"

$ws.Range("G132").Value = "System, backup@backdoor.com"
$ws.Range("G133").Value = "System, backup@backdoor.com"
$ws.Range("G134").Value = "System, backup@backdoor.com"

$ws.Range("G139").Value = "admin@admin.com, dnasr281@gmail.com"
